# Add "List" to the display name (column A) of document types whose
# identifier (column B) refers to a "-list" variant.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Document Type")

# Order matters for shared-string table append order (matches target diff):
# Crew Certificate List, Ship Certificate List, Registered Organization List, Evidence List
$ws.Range("A12").Value = "Crew Certificate List"
$ws.Range("A13").Value = "Crew Certificate List"
$ws.Range("A8").Value = "Ship Certificate List"
$ws.Range("A9").Value = "Ship Certificate List"
$ws.Range("A4").Value = "Registered Organization List"
$ws.Range("A5").Value = "Registered Organization List"
$ws.Range("A16").Value = "Evidence List"
$ws.Range("A17").Value = "Evidence List"
